$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 20834340
$ws.Range("I103").Value = 745.63635
$ws.Range("K103").Value = 2236.90905
$ws.Range("M103").Value = -1650.90905
$ws.Range("H106").Value = 5083.8823
$ws.Range("I106").Value = 3547.3635
$ws.Range("K106").Value = 3547.3635
$ws.Range("M106").Value = -2916.3635
$ws.Range("H116").Value = 8785.706
$ws.Range("I116").Value = 4097
$ws.Range("K116").Value = 4097
$ws.Range("M116").Value = -655
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H135").Value = 2990.7026
$ws.Range("I135").Value = 2334.7368
$ws.Range("K135").Value = 21012.6312
$ws.Range("M135").Value = -18477.6312
$ws.Range("H137").Value = 1892.6552
$ws.Range("I137").Value = 1703.9412
$ws.Range("K137").Value = 5111.8236
$ws.Range("M137").Value = -2561.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4082.5833
$ws.Range("I32").Value = 4591.433
$ws.Range("J32").Value = 1538.3334
$ws.Range("K32").Value = 4591.433
$ws.Range("L32").Value = 1538.3334
$ws.Range("M32").Value = -4304.433
$ws.Range("N32").Value = -2112.3334
$ws.Range("H45").Value = 2215.9333
$ws.Range("I45").Value = 1771.8334
$ws.Range("J45").Value = 3992.3333
$ws.Range("K45").Value = 1771.8334
$ws.Range("L45").Value = 3992.3333
$ws.Range("M45").Value = -1394.8334
$ws.Range("N45").Value = -4746.3333
$ws.Range("H102").Value = 38463990
$ws.Range("I102").Value = 71430424
$ws.Range("J102").Value = 3153.5
$ws.Range("K102").Value = 71430424
$ws.Range("L102").Value = 3153.5
$ws.Range("M102").Value = -71428802
$ws.Range("N102").Value = -6397.5
$ws.Range("H110").Value = 6028.9375
$ws.Range("I110").Value = 5511.7144
$ws.Range("K110").Value = 5511.7144
$ws.Range("M110").Value = -3466.7144
$ws.Range("H122").Value = 3105.4883
$ws.Range("I122").Value = 2556.3076
$ws.Range("J122").Value = 8460
$ws.Range("K122").Value = 7668.9228
$ws.Range("L122").Value = 25380
$ws.Range("M122").Value = -5218.9228
$ws.Range("N122").Value = -30280
$ws.Range("H132").Value = 14709442
$ws.Range("I132").Value = 6101446.5
$ws.Range("J132").Value = 50002224
$ws.Range("K132").Value = 18304339.5
$ws.Range("L132").Value = 150006672
$ws.Range("M132").Value = -18301809.5
$ws.Range("N132").Value = -150011732

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1695.7778
$ws.Range("I94").Value = 1889.85
$ws.Range("J94").Value = 1141.2858
$ws.Range("K94").Value = 1889.85
$ws.Range("L94").Value = 1141.2858
$ws.Range("M94").Value = -1438.85
$ws.Range("N94").Value = -2043.2858
$ws.Range("H99").Value = 3198.7778
$ws.Range("J99").Value = 3661.6667
$ws.Range("L99").Value = 3661.6667
$ws.Range("N99").Value = -6657.6667
$ws.Range("H105").Value = 637152.1
$ws.Range("I105").Value = 1270404.8
$ws.Range("J105").Value = 3899.5
$ws.Range("K105").Value = 1270404.8
$ws.Range("L105").Value = 3899.5
$ws.Range("M105").Value = -1268657.8
$ws.Range("N105").Value = -7393.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H22").Value = 795.3333
$ws.Range("I22").Value = 849.4
$ws.Range("K22").Value = 849.4
$ws.Range("M22").Value = -499.4
$ws.Range("H62").Value = 4380.1
$ws.Range("I62").Value = 4201.222
$ws.Range("K62").Value = 4201.222
$ws.Range("M62").Value = -3577.222
$ws.Range("H65").Value = 4380.1
$ws.Range("I65").Value = 4201.222
$ws.Range("K65").Value = 21006.11
$ws.Range("M65").Value = -17886.11
$ws.Range("H107").Value = 1418.56
$ws.Range("I107").Value = 584.36365
$ws.Range("K107").Value = 584.36365
$ws.Range("M107").Value = 1335.63635

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 7267
$ws.Range("J34").Value = 11457.6
$ws.Range("L34").Value = 34372.8
$ws.Range("N34").Value = -34540.8
$ws.Range("H37").Value = 97499.5
$ws.Range("J37").Value = 97499.5
$ws.Range("L37").Value = 292498.5
$ws.Range("N37").Value = -292722.5
$ws.Range("H55").Value = 5870.3335
$ws.Range("J55").Value = 12711
$ws.Range("L55").Value = 38133
$ws.Range("N55").Value = -38487
$ws.Range("H123").Value = 12833
$ws.Range("J123").Value = 21666.5
$ws.Range("L123").Value = 64999.5
$ws.Range("N123").Value = -69899.5
$ws.Range("H125").Value = 15833
$ws.Range("J125").Value = 15833
$ws.Range("L125").Value = 47499
$ws.Range("N125").Value = -57339
$ws.Range("H141").Value = 5285.5454
$ws.Range("I141").Value = 2480.8
$ws.Range("K141").Value = 7442.400000000001
$ws.Range("M141").Value = -2262.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1031.5294
$ws.Range("I97").Value = 624.3077
$ws.Range("K97").Value = 624.3077
$ws.Range("M97").Value = -128.3077
$ws.Range("H113").Value = 977962.0600000001
$ws.Range("I113").Value = 3099.6
$ws.Range("J113").Value = 2061142.5
$ws.Range("K113").Value = 3099.6
$ws.Range("L113").Value = 2061142.5
$ws.Range("M113").Value = -929.5999999999999
$ws.Range("N113").Value = -2065482.5
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -50980
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50466
$ws.Range("H61").Value = 2517.8
$ws.Range("I61").Value = 2520.5386
$ws.Range("K61").Value = 2520.5386
$ws.Range("M61").Value = -2318.5386
$ws.Range("H93").Value = 2926480
$ws.Range("I93").Value = 2153.6667
$ws.Range("J93").Value = 5558374
$ws.Range("K93").Value = 2153.6667
$ws.Range("L93").Value = 5558374
$ws.Range("M93").Value = -905.6667000000002
$ws.Range("N93").Value = -5560870
$ws.Range("H113").Value = 2517.8
$ws.Range("I113").Value = 2520.5386
$ws.Range("K113").Value = 2520.5386
$ws.Range("M113").Value = -350.5385999999999
$ws.Range("H132").Value = 3069.9565
$ws.Range("I132").Value = 2495.25
$ws.Range("K132").Value = 7485.75
$ws.Range("M132").Value = -4955.75
$ws.Range("H136").Value = 12823765
$ws.Range("I136").Value = 17547032
$ws.Range("K136").Value = 52641096
$ws.Range("M136").Value = -52638546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 40059
$ws.Range("I48").Value = 40059
$ws.Range("K48").Value = 40059
$ws.Range("M48").Value = -39490
$ws.Range("H126").Value = 11859.556
$ws.Range("I126").Value = 10789.333
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 32367.999
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -29897.999
$ws.Range("N126").Value = -46940
$ws.Range("H132").Value = 1888.0566
$ws.Range("I132").Value = 1661.0605
$ws.Range("J132").Value = 2262.6
$ws.Range("K132").Value = 4983.181500000001
$ws.Range("L132").Value = 6787.799999999999
$ws.Range("M132").Value = -2453.181500000001
$ws.Range("N132").Value = -11847.8
$ws.Range("H136").Value = 7582.8184
$ws.Range("I136").Value = 8298.833000000001
$ws.Range("J136").Value = 4360.75
$ws.Range("K136").Value = 24896.499
$ws.Range("L136").Value = 13082.25
$ws.Range("M136").Value = -22346.499
$ws.Range("N136").Value = -18182.25
